# Trade #196 closed at 2026-02-17 10:08:11 - unknown UNKNOWN +0.000%
#
# Updates:
#  - Summary sheet: refresh capital / P&L / trade-count metrics
#  - Strategy Status sheet: refresh volatility_scorer row (row 12)
#  - All Trades sheet: append the newly-closed volatility_scorer trade (#196)
#    and the newly-opened MarketMaking trade (#197)
#  - volatility_scorer sheet: append its own copy of trade #196
#  - MarketMaking sheet: append its own copy of trade #197

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(3, 2).Value = 1199.54   # Current Capital
$summary.Cells.Item(4, 2).Value = -0.46     # Total P&L $
$summary.Cells.Item(6, 2).Value = 196       # Total Trades
$summary.Cells.Item(8, 2).Value = 83        # Losing Trades
$summary.Cells.Item(9, 2).Value = 41.33     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - volatility_scorer row (row 12)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Cells.Item(12, 3).Value = 99.23     # Capital
$status.Cells.Item(12, 4).Value = 16        # Trades
$status.Cells.Item(12, 5).Value = -0.77     # P&L $
$status.Cells.Item(12, 6).Value = -0.77     # P&L %
$status.Cells.Item(12, 7).Value = 31.25     # Win Rate %

# ---------------------------------------------------------------------------
# Helper to write one trade row onto a given worksheet/row, forcing the
# date column to stay plain text instead of being auto-converted to a
# date serial number by the COM layer.
#
# NOTE: this engine's PowerShell parameter binder does not reliably bind
# named (-Param value) arguments to custom functions, so this helper uses
# plain positional parameters.
# ---------------------------------------------------------------------------
function Write-TradeRow(
    $Sheet,
    $Row,
    $TradeNum,
    $Date,
    $Time,
    $Strategy,
    $Side,
    $EntryPrice,
    $ExitPrice,
    $Status,
    $PnlPct,
    $PnlDollar,
    $CapitalAfter,
    $EntrySlippage,
    $ExitSlippage,
    $Confidence,
    $EntryReason,
    $ExitReason,
    $DurationMin
) {
    $Sheet.Cells.Item($Row, 1).Value = $TradeNum

    $Sheet.Cells.Item($Row, 2).NumberFormat = "@"
    $Sheet.Cells.Item($Row, 2).Value = $Date
    $Sheet.Cells.Item($Row, 3).Value = $Time

    $Sheet.Cells.Item($Row, 4).Value = $Strategy
    $Sheet.Cells.Item($Row, 5).Value = $Side
    $Sheet.Cells.Item($Row, 6).Value = $EntryPrice

    if ($null -ne $ExitPrice) {
        $Sheet.Cells.Item($Row, 7).Value = $ExitPrice
    }

    $Sheet.Cells.Item($Row, 8).Value = $Status
    $Sheet.Cells.Item($Row, 9).Value = $PnlPct
    $Sheet.Cells.Item($Row, 10).Value = $PnlDollar
    $Sheet.Cells.Item($Row, 11).Value = $CapitalAfter
    $Sheet.Cells.Item($Row, 12).Value = $EntrySlippage
    $Sheet.Cells.Item($Row, 13).Value = $ExitSlippage
    $Sheet.Cells.Item($Row, 14).Value = $Confidence
    $Sheet.Cells.Item($Row, 15).Value = $EntryReason

    if ($null -ne $ExitReason) {
        $Sheet.Cells.Item($Row, 16).Value = $ExitReason
    }

    $Sheet.Cells.Item($Row, 17).Value = $DurationMin
}

# ---------------------------------------------------------------------------
# All Trades sheet - append trade #196 (volatility_scorer, closed) and
# trade #197 (MarketMaking, open)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

Write-TradeRow $allTrades 197 196 `
    "2026-02-17" "10:08:04" "volatility_scorer" `
    "NEUTRAL" 0.04995 0.03 "CLOSED" `
    -39.94 -0.02 99.23 `
    0 0 0.85 `
    "Low vol market (score: inf) - ideal for market making" `
    "early_exit" 0.11

Write-TradeRow $allTrades 198 197 `
    "2026-02-17" "10:08:05" "MarketMaking" `
    "UP" 0.95 $null "OPEN" `
    0 0 100.3071991854615 `
    0 0 0.6 `
    "Normal spread capture: 19600 bps" `
    $null 0

# ---------------------------------------------------------------------------
# volatility_scorer sheet - append its own copy of trade #196
# ---------------------------------------------------------------------------
$volScorer = $wb.Worksheets.Item("volatility_scorer")

Write-TradeRow $volScorer 17 196 `
    "2026-02-17" "10:08:04" "volatility_scorer" `
    "NEUTRAL" 0.04995 0.03 "CLOSED" `
    -39.94 -0.02 99.23 `
    0 0 0.85 `
    "Low vol market (score: inf) - ideal for market making" `
    "early_exit" 0.11

# ---------------------------------------------------------------------------
# MarketMaking sheet - append its own copy of trade #197
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

Write-TradeRow $marketMaking 182 197 `
    "2026-02-17" "10:08:05" "MarketMaking" `
    "UP" 0.95 $null "OPEN" `
    0 0 100.3071991854615 `
    0 0 0.6 `
    "Normal spread capture: 19600 bps" `
    $null 0

Write-Host "Applied trade #196/#197 updates across Summary, Strategy Status, All Trades, volatility_scorer, MarketMaking"
